$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Global Reader" (row 4) loses its old remark about lockout protection.
$ws.Range("E4").ClearContents()

# 2. "Azure AD Joined Device Local Administrator" (row 9): the ANYSG-DEV-ADMINS
#    group moves from the EligiblePIM column (D) to the PermanentPIM column (C).
$groupName = $ws.Range("D9").Value2
$ws.Range("C9").Value2 = $groupName
$ws.Range("D9").Clear()

# 3. "Service Support Administrator" (row 6) gets a new remark.
$ws.Range("E6").Value = "Give to every admin to force admin MFA"

# 4. Remove the two helper "breaking.glass@alyaconsulting.ch" rows that only
#    existed to show the PIM backup account under "Privileged Role
#    Administrator" (row 5) and "Service Support Administrator" (row 7) -
#    those rows have no Role (column A) of their own.
$ws.Rows(5).Delete()
$ws.Rows(6).Delete()
